$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 71, shifting rows 71:83 down to 72:84
$ws.Rows("71").Insert()

# Populate the new row 71 with data (same structure as surrounding rows)
$ws.Range("A71").Value = 5
$ws.Range("B71").Value = "Macroferia Regional de Talca"
$ws.Range("C71").Value = "Maule"
$ws.Range("D71").Value = 44491
$ws.Range("E71").Value = 7
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100108
$ws.Range("H71").Value = "Tropicales y subtropicales"
$ws.Range("I71").Value = 100108002
$ws.Range("J71").Value = "Mango"
$ws.Range("K71").Value = "Sin especificar"
$ws.Range("L71").Value = "Primera"
$ws.Range("M71").Value = 100
$ws.Range("N71").Value = 8000
$ws.Range("O71").Value = 8000
$ws.Range("P71").Value = 8000
$ws.Range("Q71").Value = "$/bandeja 4 kilos"
$ws.Range("R71").Value = "Perú"
$ws.Range("S71").Value = 2000
$ws.Range("T71").Value = 4
